# Daily attendance processing - 2025-10-14 03:36:59
# Reorders the comma-separated "Recorded By" list in column G for every
# recorded session row (header "Recorded By" is in column G) by reversing
# the order of the entries in the list. Single-entry cells are unaffected
# since reversing a one-item list is a no-op.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ","
        if ($parts.Count -gt 1) {
            $trimmed = @()
            foreach ($p in $parts) {
                $trimmed += $p.Trim()
            }
            $reversed = @()
            for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }
            $newVal = [string]::Join(", ", $reversed)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
